# Model_TO_DO_List_and_Completion_Tracker.xlsx
# "All current models uploaded, with some measurement sheets."
#
# 1) Row 8 (Molding_Floor) now has its Measurement Sheet marked DONE.
# 2) A new row (19) is added for the "Torch" model (Furniture), combining
#    Torch_Bowl and Torch_Pillar, fully done except the texture status.
# 3) The sheet's scroll position / selection is updated to reflect where
#    the author was working when the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Measurement Sheet for Molding_Floor (row 8) is now done ---------
$cell = $ws.Range("E8")
$cell.Value = "DONE"
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4108
$cell.WrapText = $false
$cell.Font.Bold = $true
$cell.Font.Color = 255

# --- 2) New row for the "Torch" model -----------------------------------
$r = 19

$a = $ws.Range("A$r")
$a.Value = "Torch"
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4108
$a.WrapText = $false

$b = $ws.Range("B$r")
$b.Value = "Furniture"
$b.HorizontalAlignment = -4108
$b.VerticalAlignment = -4108
$b.WrapText = $false

$c = $ws.Range("C$r")
$c.Value = "Combination of Torch_Bowl and Torch_Pillar"
$c.WrapText = $true

$d = $ws.Range("D$r")
$d.Value = 1448
$d.VerticalAlignment = -4108

$e = $ws.Range("E$r")
$e.Value = "DONE"
$e.HorizontalAlignment = -4108
$e.VerticalAlignment = -4108
$e.WrapText = $false
$e.Font.Bold = $true
$e.Font.Color = 255

$f = $ws.Range("F$r")
$f.Value = "DONE"
$f.HorizontalAlignment = -4108
$f.VerticalAlignment = -4108
$f.WrapText = $false
$f.Font.Bold = $true
$f.Font.Color = 255

$g = $ws.Range("G$r")
$g.Value = "NOT STARTED"
$g.HorizontalAlignment = -4108
$g.VerticalAlignment = -4108
$g.WrapText = $true
$g.Font.Bold = $true
$g.Font.Color = 12611584

$ws.Rows.Item($r).RowHeight = 30

# --- 3) Update the view: scroll down a bit and move the selection -------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H8").Select()
